$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.913500000000006
$ws.Range("B7").Value = 4.976699999999997
$ws.Range("C7").Value = -13.9926
$ws.Range("C15").Value = -14.57119999999999
$ws.Range("B16").Value = 7.277499999999996
$ws.Range("D16").Value = -8.492100000000001
$ws.Range("D19").Value = -8.354799999999992
$ws.Range("C21").Value = -12.38610000000001
$ws.Range("C22").Value = -12.58860000000001
$ws.Range("C23").Value = -12.37410000000001
$ws.Range("B28").Value = 5.598600000000001
$ws.Range("B29").Value = 5.154700000000006
$ws.Range("B32").Value = 7.525499999999994
$ws.Range("C34").Value = -11.71960000000001
$ws.Range("E34").Value = 17.5899
$ws.Range("D36").Value = -8.521299999999995
$ws.Range("B40").Value = 8.996699999999992
$ws.Range("C43").Value = -12.52979999999999
$ws.Range("E43").Value = 17.12590000000001
$ws.Range("C45").Value = -13.62219999999999
$ws.Range("D46").Value = -8.644399999999997
$ws.Range("E48").Value = 17.50610000000001
$ws.Range("C50").Value = -13.98179999999999
$ws.Range("D50").Value = -8.044799999999999
$ws.Range("C51").Value = -12.286
$ws.Range("B52").Value = 5.093699999999999
$ws.Range("B57").Value = 5.136099999999996
$ws.Range("B66").Value = 5.935600000000001
$ws.Range("C66").Value = -11.29830000000001
$ws.Range("C67").Value = -11.27179999999999
$ws.Range("E70").Value = 17.68550000000001
$ws.Range("E73").Value = 17.53710000000001
$ws.Range("C79").Value = -11.46290000000001
$ws.Range("C84").Value = -12.7842
$ws.Range("E87").Value = 16.46
$ws.Range("C92").Value = -11.1996
$ws.Range("E92").Value = 18.32790000000001
$ws.Range("D95").Value = -8.340399999999999
$ws.Range("C97").Value = -11.6608
$ws.Range("D97").Value = -8.469599999999993
$ws.Range("B100").Value = 5.3
$ws.Range("E101").Value = 16.79800000000001
